$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '91.388.31'
$ws.Cells.Item(2, 5).Value = '  +1.77%  '

$ws.Cells.Item(3, 4).Value = '3.166.39'
$ws.Cells.Item(3, 5).Value = '  +2.81%  '

$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  +0.49%  '

$ws.Cells.Item(5, 4).Value = '''238.97'
$ws.Cells.Item(5, 5).Value = '  +1.27%  '

$ws.Cells.Item(6, 4).Value = '''619.18'
$ws.Cells.Item(6, 5).Value = '  +0.15%  '

$ws.Cells.Item(7, 5).Value = '  +6.05%  '

$ws.Cells.Item(8, 4).Value = '''0.374'
$ws.Cells.Item(8, 5).Value = '  +3.31%  '

$ws.Cells.Item(9, 4).Value = '''0.999'
$ws.Cells.Item(9, 5).Value = '  -0.12%  '

$ws.Cells.Item(10, 4).Value = '''0.741'
$ws.Cells.Item(10, 5).Value = '  +3.57%  '

$ws.Cells.Item(11, 5).Value = '  -22.42%  '

$ws.Cells.Item(12, 5).Value = '  +2.40%  '

$ws.Cells.Item(13, 4).Value = '''0.0000246'
$ws.Cells.Item(13, 5).Value = '  -1.69%  '

$ws.Cells.Item(14, 4).Value = '''35.26'
$ws.Cells.Item(14, 5).Value = '  -0.29%  '

$ws.Cells.Item(15, 4).Value = '''5.54'
$ws.Cells.Item(15, 5).Value = '  +3.12%  '

$ws.Cells.Item(16, 4).Value = '91.226.53'
$ws.Cells.Item(16, 5).Value = '  +1.68%  '

$ws.Cells.Item(17, 4).Value = '3.744.46'
$ws.Cells.Item(17, 5).Value = '  +2.03%  '

$ws.Cells.Item(18, 4).Value = '3.156.52'
$ws.Cells.Item(18, 5).Value = '  +2.09%  '

$ws.Cells.Item(19, 5).Value = '  -2.78%  '

$ws.Cells.Item(20, 4).Value = '''15.11'
$ws.Cells.Item(20, 5).Value = '  +9.30%  '

$ws.Cells.Item(21, 5).Value = '  +8.03%  '

$ws.Cells.Item(22, 4).Value = '''0.0000202'
$ws.Cells.Item(22, 5).Value = '  -4.96%  '

$ws.Cells.Item(23, 4).Value = '''442.02'
$ws.Cells.Item(23, 5).Value = '  +1.87%  '

$ws.Cells.Item(24, 4).Value = '''9.15'
$ws.Cells.Item(24, 5).Value = '  +4.00%  '

$ws.Cells.Item(25, 4).Value = '''5.74'
$ws.Cells.Item(25, 5).Value = '  -0.26%  '

$ws.Cells.Item(26, 2).Value = 'Aptos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(26, 4).Value = '''11.90'
$ws.Cells.Item(26, 5).Value = '  +1.03%  '

$ws.Cells.Item(27, 2).Value = 'Litecoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(27, 4).Value = '''82.26'
$ws.Cells.Item(27, 5).Value = '  -4.98%  '

$ws.Cells.Item(28, 4).Value = '3.321.34'
$ws.Cells.Item(28, 5).Value = '  +1.46%  '

$ws.Cells.Item(29, 5).Value = '  -0.05%  '

$ws.Cells.Item(30, 4).Value = '''0.233'
$ws.Cells.Item(30, 5).Value = '  +20.00%  '

$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(31, 4).Value = '''0.123'
$ws.Cells.Item(31, 5).Value = '  +39.23%  '

$ws.Cells.Item(32, 2).Value = 'Cronos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(32, 4).Value = '''0.171'
$ws.Cells.Item(32, 5).Value = '  +9.39%  '

$ws.Cells.Item(33, 4).Value = '''9.34'
$ws.Cells.Item(33, 5).Value = '  +2.67%  '

$ws.Cells.Item(34, 5).Value = '  +11.42%  '

$ws.Cells.Item(35, 4).Value = '''0.940'
$ws.Cells.Item(35, 5).Value = '  -5.90%  '

$ws.Cells.Item(36, 4).Value = '''7.61'
$ws.Cells.Item(36, 5).Value = '  +6.77%  '

$ws.Cells.Item(37, 4).Value = '''26.33'
$ws.Cells.Item(37, 5).Value = '  +2.70%  '

$ws.Cells.Item(38, 4).Value = '''505.76'
$ws.Cells.Item(38, 5).Value = '  +1.90%  '

$ws.Cells.Item(39, 4).Value = '''1.35'
$ws.Cells.Item(39, 5).Value = '  +6.99%  '

$ws.Cells.Item(40, 4).Value = '''1.92'
$ws.Cells.Item(40, 5).Value = '  +2.24%  '

$ws.Cells.Item(41, 4).Value = '''0.447'
$ws.Cells.Item(41, 5).Value = '  +12.65%  '

$ws.Cells.Item(42, 4).Value = '''3.82'
$ws.Cells.Item(42, 5).Value = '  +6.05%  '

$ws.Cells.Item(43, 4).Value = '''3.45'
$ws.Cells.Item(43, 5).Value = '  -7.38%  '

$ws.Cells.Item(44, 4).Value = '''22.17'
$ws.Cells.Item(44, 5).Value = '  +0.33%  '

$ws.Cells.Item(46, 4).Value = '''159.55'
$ws.Cells.Item(46, 5).Value = '  +5.17%  '

$ws.Cells.Item(47, 4).Value = '''0.710'
$ws.Cells.Item(47, 5).Value = '  +4.98%  '

$ws.Cells.Item(48, 4).Value = '''1.92'
$ws.Cells.Item(48, 5).Value = '  +3.33%  '

$ws.Cells.Item(49, 4).Value = '''1.36'
$ws.Cells.Item(49, 5).Value = '  +4.20%  '

$ws.Cells.Item(50, 4).Value = '''44.09'
$ws.Cells.Item(50, 5).Value = '  -0.62%  '

$ws.Cells.Item(51, 4).Value = '''4.42'
$ws.Cells.Item(51, 5).Value = '  +1.31%  '
